$d = $word.ActiveDocument

# --- Change 1: Replace the "User Characteristics" block with the new expanded content ---
$p4 = $d.Paragraphs.Item(4)
$p13 = $d.Paragraphs.Item(13)
$rng1 = $d.Range($p4.Range.Start, $p13.Range.End)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/></w:rPr><w:t xml:space="preserve">The purpose of this Software Requirements Specification document is to propose a system which will help the NTU staff to manage the stores. They currently do not possess a system which allows them to monitor the quantity of the products being consumed or to allow the (consumer) staffs to efficiently take the products. Instead, they write down the details on a piece of paper, from which the data is manually stored in 2 separate systems by the management staff; QuickBooks, which is responsible for creating an invoice for each department, and Sage, to update the database. The 2 systems mentioned are difficult to use, in addition to being time consuming and lack a lot of features making it hard for the staff to manage the stores. The proposed system will overcome this by combining Sage and QuickBooks along with providing additional features, such as the functionality to scan the products using external hardware, which will make it easier for the consumer staff to take items and the management staff to manage, as the system will update the database automatically.   </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/></w:rPr><w:t>The system is required to be able to; manage sock, create invoice for each department, create log of products taken by (which) staff, have unit of conversion between the unit supplied into the amount consumed, add new products on the catalogue, register products received from shipments, provide directions to help staff find items, notify the staff about important events, such as low stock, or shipment date. The system will deploy a very simple, easy to use UI with minimal input which will require no prior knowledge on using management systems. It will also provide staff with a ‘basket’ feature which will allow them to take bulk of items at once. There will also be an option for returning items, so if staff accidently took more than required, they can simply decrease the quantity on the checkout interface. The catalogue will have pictures with a brief description of the items to help consumer staff identify them easily. All the data will be backed up on an external (cloud-based) database to avoid data loss.</w:t></w:r></w:p><w:p/><w:p><w:bookmarkStart w:id="0" w:name="_Toc19884847"/><w:r><w:t>User Characteristics</w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> –</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The proposed product management system will include three kinds of users. The system privileges/features available to each of these will vary, as will their experience </w:t></w:r><w:r><w:t>and knowledge of product management. The three types of users will include:</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Consumer</w:t></w:r></w:p><w:p><w:r><w:t>The consumers will be comprised of university researchers/lectures who require the products available at the store. These users will be limited to removing and returning products to the store</w:t></w:r><w:r><w:t>. They will not require much experience with virtual systems due to their limited usage of the system and the intuitive user interface which will walk them through the product removal process. The consumers will only need to understand the basic fundamentals of computer usage and graphical user interfaces in order to utilise the system. Interface basics such as clicking a plus to increment the quantity of a product, back arrow to move to the previous page etc. will be all that is needed of this user.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Staff</w:t></w:r></w:p><w:p><w:r><w:t>The staff will be the employees working in the store that are in charge of managing stock, receiving shipments and assisting consumers. This type of user will have greater system privileges than the consumer as they will be responsible for the systems store management features. The proposed system is designed to automate a large amount of the manual work needed in the current system and as such the work load is lessened. However</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>the staff will still be required to input data into the system such as shipment codes so that the system can update the log and the stock database. This will require staff to have some experience with data input in addition to graphical user interfaces. Experience with systems such as MS Excel will transfer easily to the proposed system.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Admin</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The admin will be the store employees that have greater control over how the store is managed and run. </w:t></w:r><w:bookmarkStart w:id="1" w:name="_Toc19884848"/><w:r><w:t>The admin users will be responsible for managing the catalogue of items available to the consumers in addition to viewing logs of removed products and received shipments. Because admins will have greater control over the store as a whole, they will require management and stock control experience in order to fully utilise the system.</w:t></w:r><w:bookmarkStart w:id="2" w:name="_GoBack"/><w:bookmarkEnd w:id="2"/><w:r><w:t xml:space="preserve"> The level of technical expertise will be similar to that of the staff. </w:t></w:r></w:p><w:p><w:r><w:t>Assumptions</w:t></w:r><w:bookmarkEnd w:id="1"/><w:r><w:t xml:space="preserve"> –</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

Write-Host "Change 1 done. Paragraph count now:" $d.Paragraphs.Count

# --- Change 2: remove lastRenderedPageBreak from "Also describe any items..." paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Also describe any items that will constrain the design options*") {
        $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/></w:rPr><w:t>Also describe any items that will constrain the design options, including</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml2)
        break
    }
}

Write-Host "Change 2 done"

# --- Change 3: split "Audit functions (audit trail, log files, etc.)" run with proofErr spans ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Audit functions (audit trail*") {
        $xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Audit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/><w:lang w:val="fr-FR"/></w:rPr><w:t>functions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> (audit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/><w:lang w:val="fr-FR"/></w:rPr><w:t>trail</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Verdana" w:cs="Verdana"/><w:lang w:val="fr-FR"/></w:rPr><w:t>, log files, etc.)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml3)
        break
    }
}

Write-Host "Change 3 done"

Write-Host "Final paragraph count:" $d.Paragraphs.Count
